# Chain of Trust concept
# Applies the MSP.ORG1/MSP.ORG2 swap + repositioning of the "RCA1" cluster
# on slide 13, and removes the now-superfluous dotted connector between the
# "Rounded Rectangle 76" and "Rounded Rectangle 39" boxes.

function EmuToPt($emu) {
    # PowerPoint COM measures Left/Top/Width/Height in points (1 pt = 12700 EMU).
    # Aim for the middle of the target EMU's point-bucket so the
    # point -> EMU round-trip lands exactly back on $emu.
    return ($emu + 0.5) / 12700
}

function FindShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# --- Swap the MSP.ORG1 / MSP.ORG2 labels -----------------------------------
(FindShapeByName $s "Rounded Rectangle 49").TextFrame.TextRange.Text = "MSP.ORG2"
(FindShapeByName $s "Rounded Rectangle 57").TextFrame.TextRange.Text = "MSP.ORG1"
(FindShapeByName $s "Rounded Rectangle 76").TextFrame.TextRange.Text = "MSP.ORG2"
(FindShapeByName $s "Rounded Rectangle 39").TextFrame.TextRange.Text = "MSP.ORG1"
(FindShapeByName $s "Rounded Rectangle 65").TextFrame.TextRange.Text = "MSP.ORG1"

# --- Reposition the "RCA1" box (Rounded Rectangle 98) ----------------------
$rca1 = FindShapeByName $s "Rounded Rectangle 98"
$rca1.Left = EmuToPt 8068639
$rca1.Top = EmuToPt 3754548

# --- Reposition/resize the two elbow connectors feeding RCA1 ---------------
$elbow111 = FindShapeByName $s "Elbow Connector 111"
$elbow111.Left = EmuToPt 7353729
$elbow111.Top = EmuToPt 2752102
$elbow111.Width = EmuToPt 826294
$elbow111.Height = EmuToPt 1178598

$elbow122 = FindShapeByName $s "Elbow Connector 122"
$elbow122.Left = EmuToPt 7086939
$elbow122.Top = EmuToPt 3844188
$elbow122.Width = EmuToPt 791507
$elbow122.Height = EmuToPt 1746965

# --- Reposition the "Rounded Rectangle 39" box (now MSP.ORG1) --------------
$rr39 = FindShapeByName $s "Rounded Rectangle 39"
$rr39.Left = EmuToPt 1718030
$rr39.Top = EmuToPt 3878292

# --- Reposition/resize the long straight arrow connector below it ----------
$arrow26 = FindShapeByName $s "Straight Arrow Connector 26"
$arrow26.Left = EmuToPt 2747092
$arrow26.Top = EmuToPt 4038232
$arrow26.Width = EmuToPt 5321547
$arrow26.Height = EmuToPt 1

# --- Remove the now-unneeded dotted connector between box 76 and box 39 ----
$straight42 = FindShapeByName $s "Straight Connector 42"
if ($straight42 -ne $null) {
    $straight42.Delete()
}
